$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on cells whose new values would otherwise
# be auto-converted to numbers by Excel (they are stored as text in the source).
$textCells = @("D5", "D6", "D7", "D8", "D11", "D12", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D36", "D37", "D38", "D44", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.836.29'
$ws.Range('E2').Value = '  +4.81%  '
$ws.Range('D3').Value = '3.413.22'
$ws.Range('E3').Value = '  +3.72%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '596.03'
$ws.Range('E5').Value = '  +7.99%  '
$ws.Range('D6').Value = '188.20'
$ws.Range('E6').Value = '  +1.54%  '
$ws.Range('D7').Value = '0.603'
$ws.Range('E7').Value = '  +4.92%  '
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('E9').Value = '  +6.00%  '
$ws.Range('E10').Value = '  +3.21%  '
$ws.Range('D11').Value = '47.86'
$ws.Range('E11').Value = '  +5.14%  '
$ws.Range('D12').Value = '0.0000283'
$ws.Range('E12').Value = '  +8.19%  '
$ws.Range('D13').Value = '3.958.62'
$ws.Range('E13').Value = '  +3.42%  '
$ws.Range('D14').Value = '645.63'
$ws.Range('E14').Value = '  +11.80%  '
$ws.Range('D15').Value = '8.66'
$ws.Range('E15').Value = '  +3.01%  '
$ws.Range('D16').Value = '68.925.71'
$ws.Range('E16').Value = '  +4.86%  '
$ws.Range('D17').Value = '3.421.39'
$ws.Range('E17').Value = '  +3.19%  '
$ws.Range('E18').Value = '  +2.15%  '
$ws.Range('D19').Value = '18.18'
$ws.Range('E19').Value = '  +3.13%  '
$ws.Range('D20').Value = '11.20'
$ws.Range('E20').Value = '  +3.84%  '
$ws.Range('D21').Value = '0.918'
$ws.Range('E21').Value = '  +3.55%  '
$ws.Range('D22').Value = '18.11'
$ws.Range('E22').Value = '  +0.87%  '
$ws.Range('D23').Value = '5.13'
$ws.Range('E23').Value = '  +3.31%  '
$ws.Range('D24').Value = '100.66'
$ws.Range('E24').Value = '  +2.86%  '
$ws.Range('D25').Value = '4.11'
$ws.Range('E25').Value = '  +4.29%  '
$ws.Range('E26').Value = '  +8.07%  '
$ws.Range('D27').Value = '9.86'
$ws.Range('E27').Value = '  +5.87%  '
$ws.Range('D28').Value = '33.00'
$ws.Range('E28').Value = '  +8.94%  '
$ws.Range('D29').Value = '8.77'
$ws.Range('E29').Value = '  +5.20%  '
$ws.Range('D30').Value = '6.90'
$ws.Range('E30').Value = '  +5.04%  '
$ws.Range('D31').Value = '616.45'
$ws.Range('E31').Value = '  +8.49%  '
$ws.Range('D32').Value = '3.86'
$ws.Range('E32').Value = '  +5.28%  '
$ws.Range('D33').Value = '4.047.29'
$ws.Range('E33').Value = '  +8.77%  '
$ws.Range('D34').Value = '11.20'
$ws.Range('E34').Value = '  +3.84%  '
$ws.Range('E35').Value = '  +4.55%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = '57.07'
$ws.Range('E37').Value = '  +3.42%  '
$ws.Range('D38').Value = '2.83'
$ws.Range('E38').Value = '  +9.48%  '
$ws.Range('E39').Value = '  +5.38%  '
$ws.Range('E40').Value = '  +7.01%  '
$ws.Range('E41').Value = '  +0.87%  '
$ws.Range('D42').Value = '0.0₃0715'
$ws.Range('E42').Value = '  +4.66%  '
$ws.Range('E43').Value = '  +4.79%  '
$ws.Range('D44').Value = '3.44'
$ws.Range('E44').Value = '  +2.80%  '
$ws.Range('E45').Value = '  +5.53%  '
$ws.Range('E46').Value = '  +2.62%  '
$ws.Range('E47').Value = '  +5.34%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '1.39'
$ws.Range('E48').Value = '  +13.39%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('D50').Value = '129.97'
$ws.Range('E50').Value = '  +2.20%  '
$ws.Range('D51').Value = '7.85'
$ws.Range('E51').Value = '  +8.40%  '

# Restore default (Normal) style on the text-forced cells so no stray
# number-format style is left behind.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}